$p = $ppt.ActivePresentation
$p.Slides.Item(9).Delete()
